$wb = $excel.ActiveWorkbook

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

foreach ($ws in $wb.Worksheets) {
    $first = $ws.Cells.Find($oldTimestamp)
    if ($first -eq $null) {
        continue
    }

    $firstAddr = $first.Address()
    $cellsToFix = New-Object System.Collections.ArrayList

    $cell = $first
    do {
        [void]$cellsToFix.Add($cell.Address())
        $cell = $ws.Cells.FindNext($cell)
    } while ($cell.Address() -ne $firstAddr)

    foreach ($addr in $cellsToFix) {
        $target = $ws.Range($addr)
        $val = $target.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldTimestamp)) {
            $target.Value2 = $val.Replace($oldTimestamp, $newTimestamp)
        }
    }
}
